$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Apply formatting for the four new rows (32-35) first, by copying
# the look of existing rows, so no brand-new style entries are
# created in styles.xml (reuses the workbook's existing bordered /
# wrap-text cell formats).
# ------------------------------------------------------------------

# Rows 32 & 33 look like row 28 (border on A,B,D,E ; wrapped border on C)
$ws.Range("A28:E28").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)
$ws.Range("A28:E28").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)

# Rows 34 & 35: A,B,D,E look like row 25 ; C is a wrapped cell like C2
$ws.Range("A25:E25").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)
$ws.Range("A25:E25").Copy()
$ws.Range("A35:E35").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C35").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Fill in the new Customer Care test case data. The write order
# below matches how the rows were actually authored (rows 32 & 33
# completed in full, then rows 34 & 35 had their Jira id / Description
# columns filled before their TCID column).
# ------------------------------------------------------------------

# --- Row 32: DRAIAMCC004 ---
$ws.Cells.Item(32, 1).Value = "DRAIAMCC004"
$ws.Cells.Item(32, 2).Value = "OPQA-5229||OPQA-5230"
$ws.Cells.Item(32, 3).Value = 'Verify that Page should change header title for different title||Ensure that the page has "Support Request" and "Call us" sections.'
$ws.Cells.Item(32, 4).Value = "Y"
$ws.Cells.Item(32, 5).Value = ""

# --- Row 33: DRAIAMCC003 ---
$ws.Cells.Item(33, 1).Value = "DRAIAMCC003"
$ws.Cells.Item(33, 2).Value = "OPQA-5174"
$ws.Cells.Item(33, 3).Value = "Verify that the user should be able to select the issue type/category of the issue as an option."
$ws.Cells.Item(33, 4).Value = "Y"
$ws.Cells.Item(33, 5).Value = ""

# --- Rows 34 & 35: Jira id (B) and Description (C) filled first for both rows ---
$ws.Cells.Item(34, 2).Value = "OPQA-5169||OPQA-5170"
$ws.Cells.Item(34, 3).Value = "Verify that error messages/validation alerts `"Please enter at least 2 characters for name`" should be displayed when 'name' field is empty or 'name' field contains less than two characters.||Verify that error messages/validation alerts `"Please enter at least 2 characters for Organization Name `" should be displayed when 'Organization Name' field is empty or 'Organization Name' field contains less than two characters."

$ws.Cells.Item(35, 2).Value = "OPQA-5171||OPQA-5172"
$ws.Cells.Item(35, 3).Value = 'Verify that error messages/validation alerts " Incorrect email address format. Please try again." should be displayed when user enters incorrect email address.||Verify that error messages/validation alerts "Incorrect phone number format. Please try again.." should be displayed when user enters incorrect phone number.'

# --- TCID (column A) for rows 34 & 35 filled last ---
$ws.Cells.Item(34, 1).Value = "DRAIAMCC001"
$ws.Cells.Item(35, 1).Value = "DRAIAMCC002"

$ws.Cells.Item(34, 4).Value = "Y"
$ws.Cells.Item(34, 5).Value = ""
$ws.Cells.Item(35, 4).Value = "Y"
$ws.Cells.Item(35, 5).Value = ""

# Row heights: 34 & 35 are taller (wrapped, multi-line) rows, 45pt
$ws.Rows.Item(34).RowHeight = 45
$ws.Rows.Item(35).RowHeight = 45

# ------------------------------------------------------------------
# View state: scroll + selection to match the author's final cursor
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("F35").Select()
